$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Initial Position Single Thread"
$ws2 = $wb.Worksheets.Item(2)   # "vs other Enignes"

# ---------------------------------------------------------------------------
# 1) Restyle P92: it no longer uses the "note" highlight style, it now uses
#    the plain label style (same direct formatting as sheet2!B1/C1/D1).
# ---------------------------------------------------------------------------
$ws2.Range("B1").Copy() | Out-Null
$ws1.Range("P92").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Add three new data rows (96, 97, 98) continuing the existing table,
#    mirroring the layout of rows 92/93/94 (same per-column styling).
# ---------------------------------------------------------------------------
$ws1.Range("A92:N92").Copy() | Out-Null
$ws1.Range("A96:N96").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws1.Range("C93:N93").Copy() | Out-Null
$ws1.Range("C97:N97").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws1.Range("I94:N94").Copy() | Out-Null
$ws1.Range("I98:N98").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 96
$ws1.Range("A96").Value = 46027
$ws1.Range("C96").Value = 4
$ws1.Range("D96").Value = 206603
$ws1.Range("E96").Value = 218
$ws1.Range("F96").Formula = "=D96/E96*1000"
$ws1.Range("G96").Formula = "=(E92-E96)/E92"
$ws1.Range("H96").Formula = "=(F96-80000000)/80000000"
$ws1.Range("I96").Value = 4
$ws1.Range("J96").Value = 197281
$ws1.Range("K96").Value = 7
$ws1.Range("L96").Formula = "=J96/K96*1000"
$ws1.Range("M96").Formula = "=(K92-K96)/K92"
$ws1.Range("N96").Formula = "=(L96-80000000)/80000000"
$ws1.Range("P96").Value = "changed board representation to the classic Little-Endian Rank-File Mapping"

# Row 97
$ws1.Range("C97").Value = 5
$ws1.Range("D97").Value = 5072212
$ws1.Range("E97").Value = 5403
$ws1.Range("F97").Formula = "=D97/E97*1000"
$ws1.Range("G97").Formula = "=(E93-E97)/E93"
$ws1.Range("H97").Formula = "=(F97-80000000)/80000000"
$ws1.Range("I97").Value = 5
$ws1.Range("J97").Value = 4880523
$ws1.Range("K97").Value = 170
$ws1.Range("L97").Formula = "=J97/K97*1000"
$ws1.Range("M97").Formula = "=(K93-K97)/K93"
$ws1.Range("N97").Formula = "=(L97-80000000)/80000000"
$ws1.Range("P97").Value = "downgrade to 32GB of RAM with minor performance"

# Row 98
$ws1.Range("I98").Value = 6
$ws1.Range("J98").Value = 119060324
$ws1.Range("K98").Value = 4295
$ws1.Range("L98").Formula = "=J98/K98*1000"
$ws1.Range("M98").Formula = "=(K94-K98)/K94"
$ws1.Range("N98").Formula = "=(L98-80000000)/80000000"

# ---------------------------------------------------------------------------
# 3) Update the view: scroll down a bit and move the active selection.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 54
$win.ScrollColumn = 1
$ws1.Range("P96").Select() | Out-Null

Write-Host "edit applied"
